$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Each entry: table row (1-based, including blank spacer rows), column (1-based), new text
$cellUpdates = @(
    @(1, 1, "77÷4=19, 1"),
    @(1, 2, "78÷6=13, 0"),
    @(1, 3, "29÷9=3, 2"),
    @(1, 4, "64÷6=10, 4"),
    @(1, 5, "50÷9=5, 5"),

    @(5, 1, "21÷4=5, 1"),
    @(5, 2, "82÷6=13, 4"),
    @(5, 3, "79÷3=26, 1"),
    @(5, 4, "36÷9=4, 0"),
    @(5, 5, "67÷4=16, 3"),

    @(9, 1, "29÷4=7, 1"),
    @(9, 2, "98÷7=14, 0"),
    @(9, 3, "60÷2=30, 0"),
    @(9, 4, "44÷8=5, 4"),
    @(9, 5, "99÷8=12, 3"),

    @(13, 1, "46÷9=5, 1"),
    @(13, 2, "46÷8=5, 6"),
    @(13, 3, "82÷9=9, 1"),
    @(13, 4, "15÷6=2, 3"),
    @(13, 5, "99÷2=49, 1"),

    @(17, 1, "40÷6=6, 4"),
    @(17, 2, "31÷5=6, 1"),
    @(17, 3, "83÷9=9, 2"),
    @(17, 4, "38÷9=4, 2"),
    @(17, 5, "93÷3=31, 0")
)

foreach ($u in $cellUpdates) {
    $row = $u[0]
    $col = $u[1]
    $newText = $u[2]

    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, then set the new text.
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}
